# Applies updated price/volume figures to the cryptos price table (columns D and E).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# row => @{ D = <new price text>; E = <new volume text> }
$updates = @{
    2 = @{ D='60.938.56'; E='  +1.13%  ' }
    3 = @{ D='3.378.46'; E='  +0.03%  ' }
    4 = @{ E='  +0.03%  ' }
    5 = @{ D='570.50'; E='  -0.15%  ' }
    6 = @{ D='141.07'; E='  -0.25%  ' }
    7 = @{ E='  -0.01%  ' }
    8 = @{ D='0.474'; E='  -0.09%  ' }
    9 = @{ D='7.63'; E='  +2.11%  ' }
    10 = @{ D='0.122'; E='  -1.10%  ' }
    11 = @{ E='  -2.01%  ' }
    12 = @{ D='3.965.89'; E='  +0.29%  ' }
    13 = @{ E='  +1.86%  ' }
    14 = @{ D='27.75'; E='  -1.43%  ' }
    15 = @{ D='3.388.63'; E='  +0.21%  ' }
    16 = @{ E='  -0.37%  ' }
    17 = @{ D='61.042.66'; E='  +0.99%  ' }
    18 = @{ E='  -2.58%  ' }
    19 = @{ D='13.58'; E='  -3.57%  ' }
    20 = @{ D='8.90'; E='  -2.39%  ' }
    21 = @{ D='382.41'; E='  -1.68%  ' }
    22 = @{ D='75.51'; E='  +2.93%  ' }
    23 = @{ D='0.550'; E='  -1.80%  ' }
    24 = @{ E='  +0.47%  ' }
    25 = @{ E='  -1.46%  ' }
    26 = @{ D='3.524.26'; E='  +0.14%  ' }
    27 = @{ E='  +3.44%  ' }
    28 = @{ E='  -0.11%  ' }
    29 = @{ D='7.19'; E='  -2.77%  ' }
    30 = @{ D='7.95'; E='  -1.43%  ' }
    31 = @{ E='  -0.26%  ' }
    33 = @{ D='1.37'; E='  -3.90%  ' }
    34 = @{ D='23.14'; E='  -2.51%  ' }
    35 = @{ D='6.93'; E='  -0.11%  ' }
    36 = @{ D='165.90'; E='  -0.68%  ' }
    37 = @{ D='3.417.10'; E='  +0.29%  ' }
    38 = @{ D='4.95'; E='  -0.41%  ' }
    39 = @{ D='1.46'; E='  -2.82%  ' }
    40 = @{ D='0.0764'; E='  -1.74%  ' }
    41 = @{ D='26.69'; E='  -1.18%  ' }
    42 = @{ E='  -0.04%  ' }
    43 = @{ D='0.777'; E='  -0.54%  ' }
    44 = @{ D='4.35' }
    45 = @{ E='  -2.57%  ' }
    46 = @{ E='  -0.37%  ' }
    47 = @{ D='2.447.85'; E='  -3.21%  ' }
    48 = @{ D='22.85'; E='  -0.74%  ' }
    49 = @{ D='6.65'; E='  -2.90%  ' }
    50 = @{ D='2.13'; E='  +9.52%  ' }
    51 = @{ D='0.0261'; E='  -2.39%  ' }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    if ($rowData.ContainsKey("D")) {
        $cell = $ws.Cells.Item($row, 4)
        # Values such as "570.50" or "8.90" parse as numbers in Excel, which would
        # silently drop the significant trailing zero; force text storage first so
        # the price column keeps the exact string from the source feed.
        $priceText = $rowData["D"]
        $looksNumeric = $priceText -match "^[+-]?[0-9]*\.?[0-9]+$"
        if ($looksNumeric) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $priceText
    }
    if ($rowData.ContainsKey("E")) {
        $ws.Cells.Item($row, 5).Value = $rowData["E"]
    }
}
